{"js": "// Applies the diff: updates the date line and 25 division-problem cells\n// in the table by doing exact text search + replace, preserving formatting.\nconst replacements = [\n  [\"2025-04-25 Friday\", \"2025-04-26 Saturday\"],\n  [\"453\u00f75=90, 3\", \"360\u00f76=60, 0\"],\n  [\"297\u00f75=59, 2\", \"907\u00f73=302, 1\"],\n  [\"559\u00f76=93, 1\", \"151\u00f77=21, 4\"],\n  [\"682\u00f72=341, 0\", \"727\u00f79=80, 7\"],\n  [\"216\u00f74=54, 0\", \"876\u00f78=109, 4\"],\n  [\"359\u00f79=39, 8\", \"445\u00f73=148, 1\"],\n  [\"719\u00f75=143, 4\", \"288\u00f78=36, 0\"],\n  [\"411\u00f74=102, 3\", \"595\u00f73=198, 1\"],\n  [\"237\u00f76=39, 3\", \"603\u00f76=100, 3\"],\n  [\"925\u00f76=154, 1\", \"882\u00f74=220, 2\"],\n  [\"656\u00f73=218, 2\", \"817\u00f73=272, 1\"],\n  [\"293\u00f76=48, 5\", \"147\u00f74=36, 3\"],\n  [\"761\u00f73=253, 2\", \"121\u00f76=20, 1\"],\n  [\"817\u00f77=116, 5\", \"308\u00f73=102, 2\"],\n  [\"344\u00f72=172, 0\", \"325\u00f75=65, 0\"],\n  [\"776\u00f76=129, 2\", \"106\u00f77=15, 1\"],\n  [\"778\u00f73=259, 1\", \"901\u00f78=112, 5\"],\n  [\"163\u00f73=54, 1\", \"222\u00f79=24, 6\"],\n  [\"405\u00f73=135, 0\", \"104\u00f72=52, 0\"],\n  [\"262\u00f76=43, 4\", \"774\u00f73=258, 0\"],\n  [\"955\u00f78=119, 3\", \"938\u00f77=134, 0\"],\n  [\"803\u00f77=114, 5\", \"711\u00f72=355, 1\"],\n  [\"955\u00f79=106, 1\", \"867\u00f79=96, 3\"],\n  [\"912\u00f76=152, 0\", \"975\u00f73=325, 0\"],\n  [\"522\u00f76=87, 0\", \"534\u00f78=66, 6\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Applies the diff: updates the date line and 25 division-problem cells\n# in the table by doing exact text search + replace via Find.Execute,\n# which preserves the run's existing character formatting.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old = \"2025-04-25 Friday\"; New = \"2025-04-26 Saturday\"},\n    @{Old = \"453\u00f75=90, 3\"; New = \"360\u00f76=60, 0\"},\n    @{Old = \"297\u00f75=59, 2\"; New = \"907\u00f73=302, 1\"},\n    @{Old = \"559\u00f76=93, 1\"; New = \"151\u00f77=21, 4\"},\n    @{Old = \"682\u00f72=341, 0\"; New = \"727\u00f79=80, 7\"},\n    @{Old = \"216\u00f74=54, 0\"; New = \"876\u00f78=109, 4\"},\n    @{Old = \"359\u00f79=39, 8\"; New = \"445\u00f73=148, 1\"},\n    @{Old = \"719\u00f75=143, 4\"; New = \"288\u00f78=36, 0\"},\n    @{Old = \"411\u00f74=102, 3\"; New = \"595\u00f73=198, 1\"},\n    @{Old = \"237\u00f76=39, 3\"; New = \"603\u00f76=100, 3\"},\n    @{Old = \"925\u00f76=154, 1\"; New = \"882\u00f74=220, 2\"},\n    @{Old = \"656\u00f73=218, 2\"; New = \"817\u00f73=272, 1\"},\n    @{Old = \"293\u00f76=48, 5\"; New = \"147\u00f74=36, 3\"},\n    @{Old = \"761\u00f73=253, 2\"; New = \"121\u00f76=20, 1\"},\n    @{Old = \"817\u00f77=116, 5\"; New = \"308\u00f73=102, 2\"},\n    @{Old = \"344\u00f72=172, 0\"; New = \"325\u00f75=65, 0\"},\n    @{Old = \"776\u00f76=129, 2\"; New = \"106\u00f77=15, 1\"},\n    @{Old = \"778\u00f73=259, 1\"; New = \"901\u00f78=112, 5\"},\n    @{Old = \"163\u00f73=54, 1\"; New = \"222\u00f79=24, 6\"},\n    @{Old = \"405\u00f73=135, 0\"; New = \"104\u00f72=52, 0\"},\n    @{Old = \"262\u00f76=43, 4\"; New = \"774\u00f73=258, 0\"},\n    @{Old = \"955\u00f78=119, 3\"; New = \"938\u00f77=134, 0\"},\n    @{Old = \"803\u00f77=114, 5\"; New = \"711\u00f72=355, 1\"},\n    @{Old = \"955\u00f79=106, 1\"; New = \"867\u00f79=96, 3\"},\n    @{Old = \"912\u00f76=152, 0\"; New = \"975\u00f73=325, 0\"},\n    @{Old = \"522\u00f76=87, 0\"; New = \"534\u00f78=66, 6\"},\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $found = $range.Find.Execute(\n        $pair.Old,   # FindText\n        $false,      # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $pair.New,   # ReplaceWith\n        2            # Replace (wdReplaceAll)\n    )\n    if (-not $found) {\n        throw \"Text not found: $($pair.Old)\"\n    }\n}\n"}
